$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $ws.Range($ref).Value = "'" + $val
    $ws.Range($ref).Style = "Normal"
}

# Rename the sheet (repayment_20250915_20250915 (8) -> (9))
$ws.Name = "repayment_20250915_20250915 (9)"

# --- Row 2: Erick Ervan Dewanggga ---
$ws.Range("D2").Value = 11
Set-TextCell "E2" "5,931,659.00"
Set-TextCell "G2" "3.88"
$ws.Range("H2").Value = 667
$ws.Range("J2").Value = 2
Set-TextCell "K2" "14.95"
Set-TextCell "L2" "6.67"

# --- Row 3: Sucika Wardani ---
$ws.Range("H3").Value = 486
$ws.Range("J3").Value = 1
Set-TextCell "K3" "1.80"
Set-TextCell "L3" "3.23"

# --- Row 4: Debora Retima Sihombing ---
$ws.Range("D4").Value = 8
Set-TextCell "E4" "6,222,867.00"
Set-TextCell "G4" "3.76"
$ws.Range("H4").Value = 1.631

# --- Row 5: Fadilah Damayanti ---
$ws.Range("H5").Value = 1.1830000000000001

# --- Row 6: Axl Wicaksono ---
$ws.Range("H6").Value = 1.17

# --- Row 7: Nur Halim ---
$ws.Range("D7").Value = 9
Set-TextCell "E7" "15,512,899.00"
Set-TextCell "G7" "9.76"
$ws.Range("H7").Value = 481
$ws.Range("J7").Value = 1
Set-TextCell "K7" "34.81"
Set-TextCell "L7" "3.45"

# --- Row 8: Annisa Putri Restu ---
$ws.Range("H8").Value = 1.647

# --- Row 9: Riska Nurlita ---
$ws.Range("D9").Value = 7
Set-TextCell "E9" "884,007.00"
Set-TextCell "G9" "0.47"
$ws.Range("H9").Value = 822
Set-TextCell "K9" "0.26"

# --- Row 10: Azizah Rahmawati ---
$ws.Range("D10").Value = 8
Set-TextCell "E10" "3,042,819.00"
Set-TextCell "G10" "1.80"
$ws.Range("H10").Value = 755

# --- Row 11: Erlangga Hutama ---
$ws.Range("H11").Value = 1.2350000000000001

# --- Row 12: Ridhoi Berkat Zebua ---
$ws.Range("H12").Value = 1.329

# --- Row 13: Romli ---
$ws.Range("D13").Value = 3
Set-TextCell "E13" "722,243.00"
Set-TextCell "G13" "0.45"
$ws.Range("H13").Value = 1.7130000000000001

# --- Row 14: Aldi Taufik ---
$ws.Range("H14").Value = 780

# --- Row 15: Adistira Winditya P ---
$ws.Range("D15").Value = 1
Set-TextCell "E15" "628,614.00"
Set-TextCell "G15" "0.43"
$ws.Range("H15").Value = 651
$ws.Range("J15").Value = 1
Set-TextCell "K15" "2.35"
Set-TextCell "L15" "3.23"

# --- Row 16: Yandi Nugraha ---
$ws.Range("D16").Value = 12
Set-TextCell "E16" "2,733,557.00"
Set-TextCell "G16" "2.12"
$ws.Range("H16").Value = 1.1459999999999999
$ws.Range("J16").Value = 4
Set-TextCell "K16" "12.86"
Set-TextCell "L16" "13.33"

# --- Row 17: Wasti Feronika Sihombing ---
$ws.Range("H17").Value = 949

# --- Row 18: Nuraini ---
$ws.Range("H18").Value = 1.905
